# Add a new worksheet "Sheet2" after the existing "foo" sheet, populate it
# with a header row, and make it the active sheet/selection (matching the
# diff: new <sheet> entry, two new shared strings, and a new sheet2.xml part).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so ordering is foo, Sheet2.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "sheet2"
$ws2.Range("B1").Value = "sheet2_col2"

# Make Sheet2 the active sheet with B1 selected, as in the target workbook.
$ws2.Activate()
$ws2.Range("B1").Select()
